# Build site at 2023-04-12 14:53:07 UTC
# Fill in the LOQ4231 syllabus sheet: add the missing course-objective,
# teacher, short-syllabus, full-syllabus, method, criteria, makeup-rule and
# bibliography text, and insert a new row so the "Docentes responsaveis"
# value gets its own row before "Programa resumido".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOQ4231")

$xlPasteFormats = -4122

# --- Objetivos: (row 10) gets its course-objective text -------------------
$ws.Range("B10").Value = "Apresentar ao aluno de Engenharia conceitos básicos da Ciência Econômica"
$ws.Range("C10").Value = "Apresentar ao aluno de Engenharia conceitos básicos da Ciência Econômica"

# --- Insert a new row 13 to hold the "Docentes responsáveis" value --------
# (previously the teacher name data was misplaced under "Objetivos:" / "Método:")
$ws.Rows.Item(13).Insert()

# The inserted row copies formatting down from row 12 (column A's bold
# style); column A has no label on this row, so drop it, and pick up the
# normal/red text-column formatting for B13:C13 from a row that already
# has it.
$ws.Range("A13").Clear()
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial($xlPasteFormats)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial($xlPasteFormats)

$ws.Range("B13").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C13").Value = "5840671 - Francisco José Moreira Chaves"

# --- Programa resumido: (row 14, was row 13) short syllabus text ----------
$ws.Range("B14").Value = "História do Pensamento Econômico. Conceitos da Micro e Macroeconomia. Análise da Economia Brasileira"
$ws.Range("C14").Value = "História do Pensamento Econômico. Conceitos da Micro e Macroeconomia. Análise da Economia Brasileira"

# --- Programa: (row 16, was row 15) full syllabus text ---------------------
$programa = "1.Introdução: história do pensamento econômico." + [char]10 + "2.Microeconomia: oferta, demanda e mercado; elasticidade e estruturas de mercado (concorrência perfeita, monopólio e oligopólio)." + [char]10 + "3. Macroeconomia: teoria geral do emprego; juros e a moeda, Sistema Financeiro, Banco Central; Políticas Econômicas: inflação, crescimento, endividamento, balanço de pagamentos e comércio exterior." + [char]10 + "4.Economia brasileira"
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Método: (row 19, was row 18) teaching method text ----------------------
$ws.Range("B19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."
$ws.Range("C19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."

# --- Critério: (row 20, was row 19) grading criterion text ------------------
$ws.Range("B20").Value = "MF = (0,30*P1 + 0,60*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$ws.Range("C20").Value = "MF = (0,30*P1 + 0,60*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."

# --- Norma de recuperação: (row 21, was row 20) makeup-exam rule -----------
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."

# --- Bibliografia: (row 22, new row) bibliography text ----------------------
$biblio = "BEGG, D.; DORNBUSCH, R.; FISCHER, S. Introdução A Economia. Rio de Janeiro: Campus, 2003. " + [char]10 + "HUNT, E.K.; SHERMAN, H.J. História do Pensamento Econômico. Petrópolis: Vozes, 2000." + [char]10 + "BACHA , Edmar. Introdução à Macroeconomia: Uma perspectiva brasileira. Rio de Janeiro: Campus,1987." + [char]10 + "ROSSETTI, José Pascoal .Introdução à Economia.  9.ed. São Paulo: Atlas, 1982." + [char]10 + "SAMUELSON, P. Introdução à Economia. New York: Mc Graw-Hill Book Company."
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
